$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 9
$ws.Range("H9").Value = 133.7
$ws.Range("I9").Value = 137.44444
$ws.Range("J9").Value = 100
$ws.Range("K9").Value = 137.44444
$ws.Range("L9").Value = 100
$ws.Range("M9").Value = 31.55556000000001
$ws.Range("N9").Value = -438

# Row 31
$ws.Range("H31").Value = 1264980.5
$ws.Range("I31").Value = 1683474
$ws.Range("K31").Value = 5050422
$ws.Range("M31").Value = -5050192

# Row 34
$ws.Range("H34").Value = 18263
$ws.Range("I34").Value = 18263
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 18263
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = -18060
$ws.Range("N34").ClearContents()

# Row 36
$ws.Range("H36").Value = 18263
$ws.Range("I36").Value = 18263
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 18263
$ws.Range("L36").Value = 0
$ws.Range("M36").Value = -17548
$ws.Range("N36").ClearContents()

# Row 51
$ws.Range("H51").Value = 4000

# Row 94
$ws.Range("H94").Value = 5823.4614
$ws.Range("I94").Value = 5823.4614
$ws.Range("K94").Value = 5823.4614
$ws.Range("M94").Value = -5372.4614

# Row 112
$ws.Range("H112").Value = 2042.9286
$ws.Range("J112").Value = 2237.5417
$ws.Range("L112").Value = 6712.625100000001
$ws.Range("N112").Value = -8928.625100000001

# Row 129
$ws.Range("H129").Value = 969.0625
$ws.Range("J129").Value = 1144.4445
$ws.Range("L129").Value = 3433.3335
$ws.Range("N129").Value = -13433.3335

# Row 138
$ws.Range("H138").Value = 2252.798
$ws.Range("I138").Value = 1527.2142
$ws.Range("J138").Value = 2538.9436
$ws.Range("K138").Value = 4581.642599999999
$ws.Range("L138").Value = 7616.8308
$ws.Range("M138").Value = 558.3574000000008
$ws.Range("N138").Value = -17896.8308

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 22707.904
$ws.Range("I32").Value = 19656.537
$ws.Range("J32").Value = 34081.184
$ws.Range("K32").Value = 19656.537
$ws.Range("L32").Value = 34081.184
$ws.Range("M32").Value = -19369.537
$ws.Range("N32").Value = -34655.184

# Row 132
$ws.Range("H132").Value = 1438425.2
$ws.Range("I132").Value = 2225773
$ws.Range("J132").Value = 21199.4
$ws.Range("K132").Value = 6677319
$ws.Range("L132").Value = 63598.2
$ws.Range("M132").Value = -6674789
$ws.Range("N132").Value = -68658.20000000001

$ws = $wb.Worksheets.Item("BSM")
# Row 106
$ws.Range("H106").Value = 60000
$ws.Range("J106").Value = 60000
$ws.Range("L106").Value = 60000
$ws.Range("N106").Value = -62524

# Row 134
$ws.Range("H134").Value = 324818.88
$ws.Range("I134").Value = 446163.5
$ws.Range("J134").Value = 3612.4707
$ws.Range("K134").Value = 1338490.5
$ws.Range("L134").Value = 10837.4121
$ws.Range("M134").Value = -1335955.5
$ws.Range("N134").Value = -15907.4121

$ws = $wb.Worksheets.Item("CRP")
# Row 2
$ws.Range("H2").Value = 1669.3334
$ws.Range("I2").Value = 1669.3334
$ws.Range("K2").Value = 1669.3334
$ws.Range("M2").Value = -1556.3334

# Row 7
$ws.Range("H7").Value = 69.61539
$ws.Range("I7").Value = 35.714287
$ws.Range("K7").Value = 35.714287
$ws.Range("M7").Value = 77.285713

# Row 56
$ws.Range("H56").Value = 43333.332
$ws.Range("J56").Value = 40000
$ws.Range("L56").Value = 40000
$ws.Range("N56").Value = -41690

# Row 74
$ws.Range("H74").Value = 32725
$ws.Range("J74").Value = 32725
$ws.Range("L74").Value = 32725
$ws.Range("N74").Value = -34473

# Row 77
$ws.Range("H77").Value = 32725
$ws.Range("J77").Value = 32725
$ws.Range("L77").Value = 98175
$ws.Range("N77").Value = -106911

# Row 95
$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("N95").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
# Row 107
$ws.Range("H107").Value = 659.3
$ws.Range("I107").Value = 659.3
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 1977.9
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = -57.89999999999986
$ws.Range("N107").ClearContents()

# Row 112
$ws.Range("H112").Value = 3532
$ws.Range("I112").Value = 821.1667
$ws.Range("J112").Value = 4182.6
$ws.Range("K112").Value = 2463.5001
$ws.Range("L112").Value = 12547.8
$ws.Range("M112").Value = -1355.5001
$ws.Range("N112").Value = -14763.8

# Row 113
$ws.Range("H113").Value = 677.7895
$ws.Range("I113").Value = 550
$ws.Range("J113").Value = 692.82355
$ws.Range("K113").Value = 1650
$ws.Range("L113").Value = 2078.47065
$ws.Range("M113").Value = 520
$ws.Range("N113").Value = -6418.470649999999

# Row 131
$ws.Range("H131").Value = 15154583
$ws.Range("J131").Value = 16950848
$ws.Range("L131").Value = 50852544
$ws.Range("N131").Value = -50862624

$ws = $wb.Worksheets.Item("GSM")
# Row 5
$ws.Range("H5").Value = 33251

# Row 97
$ws.Range("H97").Value = 48540.91
$ws.Range("I97").Value = 58216.668
$ws.Range("J97").Value = 5000
$ws.Range("K97").Value = 58216.668
$ws.Range("L97").Value = 5000
$ws.Range("M97").Value = -57720.668
$ws.Range("N97").Value = -5992

# Row 98
$ws.Range("H98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("N98").ClearContents()

# Row 101
$ws.Range("H101").Value = 52891.332
$ws.Range("J101").Value = 52891.332
$ws.Range("L101").Value = 52891.332
$ws.Range("N101").Value = -59381.332

# Row 102
$ws.Range("H102").Value = 2525.6667
$ws.Range("I102").Value = 2497.5833
$ws.Range("J102").Value = 2638
$ws.Range("K102").Value = 2497.5833
$ws.Range("L102").Value = 2638
$ws.Range("M102").Value = -875.5832999999998
$ws.Range("N102").Value = -5882

# Row 109
$ws.Range("H109").Value = 30285
$ws.Range("J109").Value = 30285
$ws.Range("L109").Value = 30285
$ws.Range("N109").Value = -32365

# Row 131
$ws.Range("H131").Value = 29226
$ws.Range("J131").Value = 29226
$ws.Range("L131").Value = 29226
$ws.Range("N131").Value = -39306

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 5125
$ws.Range("I7").Value = 5833.3335
$ws.Range("J7").Value = 3000
$ws.Range("K7").Value = 5833.3335
$ws.Range("L7").Value = 3000
$ws.Range("M7").Value = -5721.3335
$ws.Range("N7").Value = -3224

# Row 126
$ws.Range("H126").Value = 5125
$ws.Range("I126").Value = 5833.3335
$ws.Range("J126").Value = 3000
$ws.Range("K126").Value = 17500.0005
$ws.Range("L126").Value = 9000
$ws.Range("M126").Value = -15030.0005
$ws.Range("N126").Value = -13940

# Row 134
$ws.Range("H134").Value = 76051.42999999999
$ws.Range("J134").Value = 76051.42999999999
$ws.Range("L134").Value = 76051.42999999999
$ws.Range("N134").Value = -86191.42999999999

# Row 135
$ws.Range("H135").Value = 123468.09
$ws.Range("J135").Value = 123468.09
$ws.Range("L135").Value = 123468.09
$ws.Range("N135").Value = -133608.09

$ws = $wb.Worksheets.Item("WVR")
# Row 4
$ws.Range("H4").Value = 3672663.2
$ws.Range("I4").Value = 11000000
$ws.Range("J4").Value = 8995
$ws.Range("K4").Value = 11000000
$ws.Range("L4").Value = 8995
$ws.Range("M4").Value = -10999887
$ws.Range("N4").Value = -9221

# Row 62
$ws.Range("H62").Value = 4300
$ws.Range("I62").Value = 4000
$ws.Range("J62").Value = 4750
$ws.Range("K62").Value = 4000
$ws.Range("L62").Value = 4750
$ws.Range("M62").Value = -3376
$ws.Range("N62").Value = -5998

# Row 65
$ws.Range("H65").Value = 4300
$ws.Range("I65").Value = 4000
$ws.Range("J65").Value = 4750
$ws.Range("K65").Value = 20000
$ws.Range("L65").Value = 23750
$ws.Range("M65").Value = -16880
$ws.Range("N65").Value = -29990

# Row 105
$ws.Range("H105").Value = 35000
$ws.Range("J105").Value = 35000
$ws.Range("L105").Value = 35000
$ws.Range("N105").Value = -41988

# Row 132
$ws.Range("H132").Value = 1973.4722
$ws.Range("I132").Value = 1268.0333
$ws.Range("J132").Value = 5500.6665
$ws.Range("K132").Value = 3804.0999
$ws.Range("L132").Value = 16501.9995
$ws.Range("M132").Value = -1274.0999
$ws.Range("N132").Value = -21561.9995

# Row 136
$ws.Range("H136").Value = 1487.675
$ws.Range("I136").Value = 1298.2858
$ws.Range("J136").Value = 1929.5834
$ws.Range("K136").Value = 3894.8574
$ws.Range("L136").Value = 5788.7502
$ws.Range("M136").Value = -1344.8574
$ws.Range("N136").Value = -10888.7502
